$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data loading: fill rows 2-12, columns A-F with incrementing values 1-10,
# with the last row (12) repeating the final value (10).
$values = @(1,2,3,4,5,6,7,8,9,10,10)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $val = $values[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $val
    }
}

$ws.Range("D4").Select()
